$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# A new product ("فازلين بيور وسط") was sold and needs a new data row
# inserted right above the totals row (old row 19 = totals, old row 20 =
# footer). After insertion: new row 19 = product #13, row 20 = totals
# (grand total bumped by the new line's price), row 21 = footer (with the
# report-generation timestamp updated).
# ---------------------------------------------------------------------------

# 1) Insert a blank row above the totals row, pushing totals/footer down.
$ws.Rows.Item(19).Insert()

# 2) Recreate the thin light-grey bottom border used by every item row
#    (this also happens to be the exact style combination already used by
#    rows 7-18, so Excel's style de-duplication reuses those style ids).
$rowRng = $ws.Range("A19:Q19")
$rowRng.Borders.Item(9).LineStyle = 1
$rowRng.Borders.Item(9).Color = 13882323

# 3) Fill in the new item's data, matching the columns used by every other
#    item row:
#      A      = sequence number
#      C:G    = item name
#      H:K    = "sold:remaining" ratio
#      L:M    = order limit (stored as text, like the rest of the column)
#      N:O    = price
#      P      = sell price (stored as text even though the column's display
#               format is numeric - matches the rest of the column)
#      Q      = number of transactions ratio

$ws.Range("A19").Value = 13
$ws.Range("C19").Value = "فازلين بيور وسط"
$ws.Range("H19").Value = "5:0"

# Column L's number format ("#,##0.##..."./";0") would otherwise coerce a
# plain "0" into a real number, so round-trip through Text format to keep it
# as a literal string like the other rows.
$lFmt = $ws.Range("L19").NumberFormat()
$ws.Range("L19").NumberFormat = "@"
$ws.Range("L19").Value = "0"
$ws.Range("L19").NumberFormat = $lFmt

$ws.Range("N19").Value = "30.00"

# Column P displays with a numeric format ("0.00") but the value itself is
# stored as text across the whole sheet - same Text round-trip trick.
$pFmt = $ws.Range("P19").NumberFormat()
$ws.Range("P19").NumberFormat = "@"
$ws.Range("P19").Value = "30.0000"
$ws.Range("P19").NumberFormat = $pFmt

$ws.Range("Q19").Value = "1:0"

# 4) Merge the new row's cells exactly like the rows above it.
$ws.Range("A19:B19").Merge()
$ws.Range("C19:G19").Merge()
$ws.Range("H19:K19").Merge()
$ws.Range("L19:M19").Merge()
$ws.Range("N19:O19").Merge()

# 5) Update the grand total (old row 19, now row 20) to include the new
#    line's sell price.
$ws.Range("P20").Value = 393.06

# 6) Update the footer's generated-at timestamp (old row 20, now row 21).
$ws.Range("A21").Value = "Wednesday, 16 July, 2025 10:28 AM"
